{"js": "// Word Office.js edit script\n// Applies three changes to the \"Deep Learning Project\" document:\n//   1. Splits the \"First, for our model...\" paragraph into five runs,\n//      moving \"First, \" later in the sentence and appending a new\n//      sentence about word embedding.\n//   2. Appends two new sentences (\"We use Keras Functional API...\")\n//      after the \"...seq2seq model as our training model. \" run.\n//   3. Inserts a new empty ListParagraph-styled paragraph right before\n//      the \"Model Evaluation\" heading.\n\nconst PKG_OPEN =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n  'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>';\nconst PKG_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrapPara(innerPXml) {\n  return PKG_OPEN + innerPXml + PKG_CLOSE;\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three target paragraphs by their (unique) current text.\nlet pOneHot = null;\nlet pSeq2Seq = null;\nlet pPredict = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"First, for our model to understand\") !== -1) {\n    pOneHot = paragraphs.items[i];\n  } else if (t.indexOf(\"seq2seq model as our training model\") !== -1) {\n    pSeq2Seq = paragraphs.items[i];\n  } else if (t.indexOf(\"probability of likelihood of that word to occur\") !== -1) {\n    pPredict = paragraphs.items[i];\n  }\n}\n\nif (!pOneHot || !pSeq2Seq || !pPredict) {\n  throw new Error(\"Could not locate one or more target paragraphs.\");\n}\n\n// --- Change 1: rewrite the \"First, for our model...\" paragraph --------\n// Replace only the paragraph's content (not the paragraph mark), keeping\n// its pPr intact, with five runs that reorder \"First,\" and append the\n// new \"However, ...\" sentence.\nconst oneHotRange = pOneHot.getRange(\"Content\");\nconst oneHotXml =\n  \"<w:p w14:paraId=\\\"6797E8C5\\\" w14:textId=\\\"4775F80A\\\" w:rsidR=\\\"00D82ADB\\\" w:rsidRDefault=\\\"00D82ADB\\\" w:rsidP=\\\"006F14F6\\\">\" +\n  \"<w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr><w:t>F</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\">or our model to understand out human and robot language, we need to do Natural language processing (NLP). </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr><w:t>First, w</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr>\" +\n  \"<w:t>e use one-hot vectors for encoder input, decoder input and decoder output.</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\"> However, considering one-hot encoding will ignore the inner meaning of the word in a sentence, \" +\n  \"so next, we use word embedding to help improving model performance.</w:t></w:r>\" +\n  \"</w:p>\";\noneHotRange.insertOoxml(wrapPara(oneHotXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Change 2: append \"We use Keras Functional API...\" sentence -------\n// Insert new runs at the end of the paragraph (after the existing,\n// untouched \"We selected \" / \"LSTM \" / \"seq2seq model...\" runs).\nconst seq2seqRange = pSeq2Seq.getRange(\"Content\");\nconst kerasXml =\n  \"<w:p>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">We use </w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr><w:t>Keras</w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Functional API as our model structure.</w:t></w:r>\" +\n  \"</w:p>\";\nseq2seqRange.insertOoxml(wrapPara(kerasXml), Word.InsertLocation.end);\nawait context.sync();\n\n// --- Change 3: insert a new empty paragraph before \"Model Evaluation\" --\nconst nextPara = pPredict.getNext();\nconst insertionPoint = nextPara.getRange(\"Start\");\nconst emptyParaXml =\n  \"<w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:rPr><w:color w:val=\\\"000000\\\" w:themeColor=\\\"text1\\\"/></w:rPr></w:pPr></w:p>\";\ninsertionPoint.insertOoxml(wrapPara(emptyParaXml), Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Word COM-interop edit script\n# Applies three changes to the \"Deep Learning Project\" document:\n#   1. Splits the \"First, for our model...\" paragraph into five runs,\n#      moving \"First, \" later in the sentence and appending a new\n#      sentence about word embedding.\n#   2. Appends two new sentences (\"We use Keras Functional API...\")\n#      after the \"...seq2seq model as our training model. \" run.\n#   3. Inserts a new empty ListParagraph-styled paragraph right before\n#      the \"Model Evaluation\" heading.\n\n$d = $word.ActiveDocument\n\n$PKG_OPEN = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>'\n$PKG_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Locate the three target paragraphs by their (unique) current text.\n$pOneHot = $null\n$pSeq2Seq = $null\n$pPredict = $null\n$pModelEval = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*First, for our model to understand*\") {\n        $pOneHot = $p\n    } elseif ($t -like \"*seq2seq model as our training model*\") {\n        $pSeq2Seq = $p\n    } elseif ($t -like \"*probability of likelihood of that word to occur*\") {\n        $pPredict = $p\n    } elseif ($t -like \"Model Evaluation*\") {\n        $pModelEval = $p\n    }\n}\n\n# --- Change 1: rewrite the \"First, for our model...\" paragraph --------\n# Replace the whole paragraph (content + mark) with one that keeps the\n# same pPr/paragraph attributes but splits the text into five runs that\n# reorder \"First,\" and append the new \"However, ...\" sentence.\n$oneHotXml = $PKG_OPEN + `\n  '<w:p w14:paraId=\"6797E8C5\" w14:textId=\"4775F80A\" w:rsidR=\"00D82ADB\" w:rsidRDefault=\"00D82ADB\" w:rsidP=\"006F14F6\">' + `\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr></w:pPr>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t>F</w:t></w:r>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr>' + `\n  '<w:t xml:space=\"preserve\">or our model to understand out human and robot language, we need to do Natural language processing (NLP). </w:t></w:r>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t>First, w</w:t></w:r>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr>' + `\n  '<w:t>e use one-hot vectors for encoder input, decoder input and decoder output.</w:t></w:r>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr>' + `\n  '<w:t xml:space=\"preserve\"> However, considering one-hot encoding will ignore the inner meaning of the word in a sentence, so next, we use word embedding to help improving model performance.</w:t></w:r>' + `\n  '</w:p>' + $PKG_CLOSE\n$pOneHot.Range.InsertXML($oneHotXml)\n\n# --- Change 2: append \"We use Keras Functional API...\" sentence -------\n# Replace the whole paragraph with the original (untouched) runs plus\n# three new trailing runs, preserving the original runs' rsidR etc.\n$seq2seqXml = $PKG_OPEN + `\n  '<w:p w14:paraId=\"47305BDD\" w14:textId=\"3B5B26BE\" w:rsidR=\"00CA34DA\" w:rsidRPr=\"00D82ADB\" w:rsidRDefault=\"00CA34DA\" w:rsidP=\"00D82ADB\">' + `\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr></w:pPr>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\">We selected </w:t></w:r>' + `\n  '<w:r w:rsidR=\"00D82ADB\"><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\">LSTM </w:t></w:r>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\">seq2seq model as our training model. </w:t></w:r>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\">We use </w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellStart\"/>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t>Keras</w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellEnd\"/>' + `\n  '<w:r><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\"> Functional API as our model structure.</w:t></w:r>' + `\n  '</w:p>' + $PKG_CLOSE\n$pSeq2Seq.Range.InsertXML($seq2seqXml)\n\n# --- Change 3: insert a new empty paragraph before \"Model Evaluation\" --\n$insertPoint = $d.Range($pModelEval.Range.Start, $pModelEval.Range.Start)\n$emptyParaXml = $PKG_OPEN + `\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:rPr><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr></w:pPr></w:p>' + $PKG_CLOSE\n$insertPoint.InsertXML($emptyParaXml)\n"}
